$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix swapped month/year columns (A and B) ---
# Originally column A held the year (2020) and column B held the month
# number (1..12), even though the headers say A=month, B=year. Swap the
# values so the data matches the headers.
for ($r = 2; $r -le 13; $r++) {
    $month = $r - 1
    $ws.Cells.Item($r, 1).Value = $month
    $ws.Cells.Item($r, 2).Value = 2020
}

# --- Add new grade_* / frequency header columns ---
$ws.Cells.Item(1, 8).Value  = "grade_total"
$ws.Cells.Item(1, 9).Value  = "grade_distance"
$ws.Cells.Item(1, 10).Value = "grade_visitation"
$ws.Cells.Item(1, 11).Value = "grade_encounters"
$ws.Cells.Item(1, 12).Value = "NEVER"
$ws.Cells.Item(1, 13).Value = "RARELY"
$ws.Cells.Item(1, 14).Value = "SOMETIMES"
$ws.Cells.Item(1, 15).Value = "FREQUENTLY"
$ws.Cells.Item(1, 16).Value = "ALWAYS"

# --- Populate the new columns for every data row ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value  = 1
    $ws.Cells.Item($r, 9).Value  = 0
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 3
    $ws.Cells.Item($r, 12).Value = 1.066
    $ws.Cells.Item($r, 13).Value = 1.08
    $ws.Cells.Item($r, 14).Value = 1.126
    $ws.Cells.Item($r, 15).Value = 1.194
    $ws.Cells.Item($r, 16).Value = 1.534
}
